$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "The study did not achieve statistical significance on the primary endpoint of change from baseline in the 17-item Hamilton Depression Rating Scale total score at Day 15, or on any secondary endpoints."
$ws.Range("N2").Value = "Did not meet any secondary endpoints"
$ws.Range("W2").Value = "Efficacy and safety"
$ws.Range("B3").Value = "0001689548-22-000111"
$ws.Range("J3").Value = "In April 2022, the Investigational New Drug application for the study of PRAX-222 was placed on clinical hold by the FDA. The company later submitted additional documentation from a toxicology study to address the hold."
$ws.Range("M3").Value = "Clinical Hold"
$ws.Range("S3").Value = "IND"
$ws.Range("X3").Value = "Completed"
$ws.Range("F4").Value = "Essential3 clinical program (Study 1 and Study 2)"
$ws.Range("J4").Value = "Based on a pre-planned interim analysis, the IDMC recommended that Study 1 be stopped for futility, as it was unlikely to meet the primary efficacy endpoint. The company has decided to continue both Study 1 and Study 2 to completion."
$ws.Range("S4").Value = "NDA"
$ws.Range("B5").Value = "0001689548-25-000058"
$ws.Range("D5").Value = "Epilepsy"
$ws.Range("F5").Value = "Photo-Paroxysmal Response (PPR) study"
$ws.Range("J5").Value = "Announced positive results from the PPR study, which is part of the ENERGY program."
$ws.Range("J6").Value = "The RADIANT study is an open label eight-week study in patients with focal onset seizures or generalized epilepsy and is part of the ENERGY program."
$ws.Range("J7").Value = "The POWER 1 study is a 12-week study in focal onset seizures and is part of the ENERGY program."
$ws.Range("O7").Value = "double-blind, placebo-controlled"
$ws.Range("Q7").Value = "Placebo"
$ws.Range("D8").Value = "Epilepsy"
$ws.Range("I8").Value = "2025H2"
$ws.Range("J8").Value = "The POWER 2 study is the third efficacy study in the ENERGY program, with enrollment planned to begin in the second half of 2025."
$ws.Range("B9").Value = "0001689548-25-000058"
$ws.Range("F9").Value = "EMBOLD study (first cohort)"
$ws.Range("F10").Value = "EMBOLD study (second cohort)"
$ws.Range("D11").Value = "developmental and epileptic encephalopathies (DEE)"
$ws.Range("F11").Value = "EMERALD study"
$ws.Range("J11").Value = "Plan to initiate the EMERALD study in a broader developmental and epileptic encephalopathies (DEE) patient population in mid-2025."
$ws.Range("B12").Value = "0001689548-25-000058"
$ws.Range("F12").Value = "EMBRAVE study (second cohort)"
$ws.Range("H12").Value = "Expected"
$ws.Range("I12").Value = "2026H1"
$ws.Range("J12").Value = "Currently enrolling the second cohort of the EMBRAVE study in Brazil, with topline results expected in the first half of 2026. Results from Part 1 were shared in the fourth quarter of 2023."
$ws.Range("R12").Value = "Brazil"
$ws.Range("W12").Value = "Topline"
$ws.Range("X12").Value = "Enrolling"

$ws.Rows.Item(13).Delete()

